# Agregando soporte para apoyos inclinados
#
# The "restric" sheet (nodo / dirección / desplazamiento / rotación) drives
# the boundary-condition table for the truss solver. Supports on direction 3
# (rollers, styled rows 5-11) used to assume a fixed rotation of 0; now that
# inclined supports are supported, that angle is no longer a known constant,
# so those cells become "not available yet" (#N/A) placeholders driven by a
# formula instead of a hard-coded 0.
#
# A small new work area is also being started below the table (around H13)
# for the upcoming inclined-support input/legend, and the selection is left
# where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("restric")
$ws.Select()

$ws.Range("D5").Formula = "=NA()"
$ws.Range("D6").Formula = "=NA()"
$ws.Range("D7").Formula = "=NA()"
$ws.Range("D8").Formula = "=NA()"
$ws.Range("D9").Formula = "=NA()"
$ws.Range("D10").Formula = "=NA()"
$ws.Range("D11").Formula = "=NA()"

# Start of the new "apoyos inclinados" work area: touch H13 with the same
# look as the other table headers (bold, centered) without giving it a
# value yet.
$ws.Range("A1").Copy()
$ws.Range("H13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the cursor where the user was last working.
$ws.Range("E12").Select()
